$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.UsedRange.Clear()

$ws.Cells.Item(1,1).Value = 'venue'
$ws.Cells.Item(1,2).Value = 'date'
$ws.Cells.Item(1,3).Value = 'result'
$ws.Cells.Item(1,4).Value = 'ownTeam'
$ws.Cells.Item(1,5).Value = 'oppTeam'
$ws.Cells.Item(1,6).Value = 'batsman'
$ws.Cells.Item(1,7).Value = 'totalRuns'
$ws.Cells.Item(1,8).Value = 'totalBalls'
$ws.Cells.Item(1,9).Value = 'total4s'
$ws.Cells.Item(1,10).Value = 'total6s'
$ws.Cells.Item(1,11).Value = 'sr'

$ws.Cells.Item(2,1).Value = ' Sharjah'
$ws.Cells.Item(2,2).Value = ' September 27 2020'
$ws.Cells.Item(2,3).Value = 'Royals won by 4 wickets (with 3 balls remaining)'
$ws.Cells.Item(2,4).Value = 'Rajasthan Royals'
$ws.Cells.Item(2,5).Value = 'Kings XI Punjab'
$ws.Cells.Item(2,6).Value = 'Tom Curran '
$ws.Cells.Item(2,7).Value = '''4'
$ws.Cells.Item(2,8).Value = '''1'
$ws.Cells.Item(2,9).Value = '''1'
$ws.Cells.Item(2,10).Value = '''0'
$ws.Cells.Item(2,11).Value = '''400.00'

$ws.Cells.Item(3,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(3,2).Value = ' September 30 2020'
$ws.Cells.Item(3,3).Value = 'KKR won by 37 runs'
$ws.Cells.Item(3,4).Value = 'Rajasthan Royals'
$ws.Cells.Item(3,5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(3,6).Value = 'Tom Curran '
$ws.Cells.Item(3,7).Value = '''54'
$ws.Cells.Item(3,8).Value = '''36'
$ws.Cells.Item(3,9).Value = '''2'
$ws.Cells.Item(3,10).Value = '''3'
$ws.Cells.Item(3,11).Value = '''150.00'

$ws.Cells.Item(4,1).Value = ' Sharjah'
$ws.Cells.Item(4,2).Value = ' September 22 2020'
$ws.Cells.Item(4,3).Value = 'Royals won by 16 runs'
$ws.Cells.Item(4,4).Value = 'Rajasthan Royals'
$ws.Cells.Item(4,5).Value = 'Chennai Super Kings'
$ws.Cells.Item(4,6).Value = 'Tom Curran '
$ws.Cells.Item(4,7).Value = '''10'
$ws.Cells.Item(4,8).Value = '''9'
$ws.Cells.Item(4,9).Value = '''1'
$ws.Cells.Item(4,10).Value = '''0'
$ws.Cells.Item(4,11).Value = '''111.11'

$ws.Cells.Item(5,1).Value = ' Abu Dhabi'
$ws.Cells.Item(5,2).Value = ' October 06 2020'
$ws.Cells.Item(5,3).Value = 'Mumbai won by 57 runs'
$ws.Cells.Item(5,4).Value = 'Rajasthan Royals'
$ws.Cells.Item(5,5).Value = 'Mumbai Indians'
$ws.Cells.Item(5,6).Value = 'Tom Curran '
$ws.Cells.Item(5,7).Value = '''15'
$ws.Cells.Item(5,8).Value = '''16'
$ws.Cells.Item(5,9).Value = '''1'
$ws.Cells.Item(5,10).Value = '''0'
$ws.Cells.Item(5,11).Value = '''93.75'
